# fix: layouts and role management try #1
#
# Updates the two "academicos" demo rows:
#   - row 1: name/email/id for the first user
#   - row 2: name/email/id for the second user
# and moves the active selection from B1 to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 -----------------------------------------------------------
$ws.Range("A1").Value = "Jeremiahs silva"
$ws.Range("B1").Value = "teste3@gmail.com"
$ws.Range("D1").Value = 100033659

# --- Row 2 -----------------------------------------------------------
$ws.Range("A2").Value = "Vilinda bastos"
$ws.Range("B2").Value = "teste4@gmail.com"
$ws.Range("D2").Value = 100023887

# --- Hyperlinks --------------------------------------------------------
# The email cells carry mailto: hyperlinks whose address + displayed text
# must follow the new email addresses above. Rebuild them from scratch
# (there's no reliable in-place "rename" on existing hyperlink entries).
$hyperlinks = $ws.Hyperlinks
$hyperlinks.Delete()
$hyperlinks.Add($ws.Range("B1"), "mailto:teste3@gmail.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "teste3@gmail.com") | Out-Null
$hyperlinks.Add($ws.Range("B2"), "mailto:teste4@gmail.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "teste4@gmail.com") | Out-Null

# Adding a hyperlink re-applies Excel's builtin "Hyperlink" cell style
# (blue + underline). Restore the original plain font so the cell
# formatting stays as it was before the edit.
$ws.Range("B1:B2").Font.Name = "Arial"
$ws.Range("B1:B2").Font.Size = 10
$ws.Range("B1:B2").Font.Underline = $false
$ws.Range("B1:B2").Font.Color = 16711680

# --- Selection ---------------------------------------------------------
$ws.Range("B2").Select() | Out-Null
